$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
for ($r = 2; $r -le 202; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $v = $cell.Value2
    $cell.Value2 = [Math]::Round($v, 0)
}
